# -----------------------------------------------------------------------
# Ebates-cnn-results.xlsx : rebuild the "micro f1" table with a second
# (held-out) run of scores, add the per-bucket lookup table in J:K,
# scale it through L6/M-column formulas, and tidy up the duplicate
# hyperlink cell-styles that had accumulated in styles.xml.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the old helper columns (I = old scaled value, L = old copy) ---
$ws.Range("I4:I14").ClearContents()
$ws.Range("L5").ClearContents()
$ws.Range("L7:L14").ClearContents()

# --- F column: new primary score; G column: the old score, kept alongside ---
$ws.Range("F5").Value = 0.83915723873000003
$ws.Range("G5").Value = 0.829531565493

$ws.Range("F6").Value = 0.878184364587
$ws.Range("G6").Value = 0.86462716881900004

$ws.Range("F7").Value = 0.87802601849600004
$ws.Range("G7").Value = 0.89717501438500002
$ws.Range("H7").Value = "Nan"

$ws.Range("F8").Value = 0.77959223072799999
$ws.Range("G8").Value = 0.76732054123000004

$ws.Range("F9").Value = 0.93806871401299996
$ws.Range("G9").Value = 0.93094244665600001

$ws.Range("F10").Value = 0.90401965502399995

$ws.Range("F11").Value = 0.92605434475000004
$ws.Range("G11").Value = 0.92478063968299995
$ws.Range("H11").Value = "done"

$ws.Range("F12").Value = 0.880859031189
$ws.Range("G12").Value = 0.87941877770499999

$ws.Range("F13").Value = 0.90432273126399998
$ws.Range("G13").Value = 0.89932004274400001

$ws.Range("F14").Value = 0.904975038733
$ws.Range("G14").Value = 0.89103115854699999

# --- J/K lookup table (bucket size -> score) + L6 scale factor ---
$ws.Range("L6").Value = 100

$ws.Range("J6").Value = 20000
$ws.Range("K6").Value = 0.83915723873000003

$ws.Range("J7").Value = 30000
$ws.Range("K7").Value = 0.86308505690000004

$ws.Range("J8").Value = 40000
$ws.Range("K8").Value = 0.878184364587

$ws.Range("J9").Value = 70000
$ws.Range("K9").Value = 0.87802601849600004

$ws.Range("J10").Value = 80000
$ws.Range("K10").Value = 0.87657635327100003

$ws.Range("J11").Value = 90000
$ws.Range("K11").Value = 0.77959223072799999

$ws.Range("J12").Value = 100000
$ws.Range("K12").Value = 0.82763087984499994

$ws.Range("J13").Value = 110000
$ws.Range("K13").Value = 0.93806871401299996

$ws.Range("J14").Value = 120000
$ws.Range("K14").Value = 0.90401965502399995

$ws.Range("J15").Value = 130000
$ws.Range("K15").Value = 0.93128839386399997

$ws.Range("J16").Value = 140000
$ws.Range("K16").Value = 0.93529972408299999

$ws.Range("J17").Value = 180000
$ws.Range("K17").Value = 0.92605434475000004

$ws.Range("J18").Value = 210000
$ws.Range("K18").Value = 0.880859031189

$ws.Range("J19").Value = 250000
$ws.Range("K19").Value = 0.90432273126399998

$ws.Range("J20").Value = 260000
$ws.Range("K20").Value = 0.904975038733

# --- M column: scaled score, M6 stand-alone, M7:M20 one shared formula ---
$ws.Range("M6").Formula = '=K6*$L$6'
$ws.Range("M7:M20").Formula = '=K7*$L$6'

# --- mirror the same scores for the small lookup table in rows 22-26 ---
$ws.Range("F22").Value = 0.86308505690000004
$ws.Range("F23").Value = 0.87657635327100003
$ws.Range("F24").Value = 0.82763087984499994
$ws.Range("F25").Value = 0.93128839386399997
$ws.Range("F26").Value = 0.93529972408299999

# --- keep the worksheet's recorded sort state in sync with the new table ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("J6"))
$sortObj.SetRange($ws.Range("J6:K20"))
$sortObj.Header = 0
$sortObj.Apply()

# --- selection left where the new scaled-score column was filled in ---
$ws.Range("M6:M20").Select()

# --- drop the duplicate "Followed Hyperlink" / "Hyperlink" cell styles ---
# (23 -> 11: keep the first 5 of each plus Normal, same as the cleaned file)
for ($i = 22; $i -ge 17; $i--) {
    $wb.Styles.Item($i).Delete()
}
for ($i = 11; $i -ge 6; $i--) {
    $wb.Styles.Item($i).Delete()
}
